# "updated for 29 tarikh"
# Update the October 2021 tracker sheet with the 29th's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revised totals for rows 3 and 5 (column B)
$ws.Range("B3").Value = 4888
$ws.Range("B5").Value = 3292

# New spend entries for the 28th/29th (column F)
$ws.Range("F30").Value = 1620
$ws.Range("F31").Value = 210

# Daily tally entries for the 29th (row 31, columns K:T)
$ws.Range("K31").Value = 2
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 1
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 2
$ws.Range("P31").Value = 2
$ws.Range("Q31").Value = 0
$ws.Range("R31").Value = 2
$ws.Range("S31").Value = 2
$ws.Range("T31").Value = 2

# Daily tally entries for the 30th (row 32, columns K:T)
$ws.Range("K32").Value = 2
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 2
$ws.Range("P32").Value = 2
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = 2
$ws.Range("S32").Value = 2
$ws.Range("T32").Value = 2

# Move the active selection to reflect where the author left off editing
$ws.Range("C25").Select()
